$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44489
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 11500
$ws.Range("O3").Value = "Provincia del Elquí"
$ws.Range("P3").Value = 288

# Row 4
$ws.Range("D4").Value = 44468
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 13000
$ws.Range("M4").Value = 12500
$ws.Range("P4").Value = 312

# Row 5
$ws.Range("D5").Value = 44446
$ws.Range("J5").Value = 160
$ws.Range("K5").Value = 12500
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 12750
$ws.Range("P5").Value = 319

# Row 6
$ws.Range("D6").Value = 44505

# Row 7
$ws.Range("D7").Value = 44494
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 11000
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 11500
$ws.Range("P7").Value = 288

# Row 8
$ws.Range("D8").Value = 44426
$ws.Range("K8").Value = 13000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 13500
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 338

# Row 9
$ws.Range("D9").Value = 44495
$ws.Range("K9").Value = 11000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 11500
$ws.Range("P9").Value = 288

# Row 10
$ws.Range("D10").Value = 44454
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 14000
$ws.Range("M10").Value = 13500
$ws.Range("P10").Value = 338

# Row 11
$ws.Range("D11").Value = 44417
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 16000
$ws.Range("M11").Value = 15500
$ws.Range("P11").Value = 388

# Row 12
$ws.Range("D12").Value = 44512
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 11000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 11500
$ws.Range("P12").Value = 288

# Row 13
$ws.Range("D13").Value = 44498
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 10500
$ws.Range("L13").Value = 11000
$ws.Range("M13").Value = 10750
$ws.Range("P13").Value = 269

# Row 14
$ws.Range("D14").Value = 44484

# Row 15
$ws.Range("D15").Value = 44510

# Row 16
$ws.Range("D16").Value = 44435
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 14500
$ws.Range("P16").Value = 362

# Row 17
$ws.Range("D17").Value = 44516
$ws.Range("K17").Value = 11000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 11500
$ws.Range("P17").Value = 288

# Row 18
$ws.Range("D18").Value = 44503
$ws.Range("O18").Value = "Provincia del Elquí"

# Row 19
$ws.Range("D19").Value = 44420
$ws.Range("K19").Value = 13000
$ws.Range("L19").Value = 14000
$ws.Range("M19").Value = 13500
$ws.Range("P19").Value = 338

# Row 20
$ws.Range("D20").Value = 44427
$ws.Range("K20").Value = 13000
$ws.Range("L20").Value = 14000
$ws.Range("M20").Value = 13500
$ws.Range("P20").Value = 338

# Row 21
$ws.Range("D21").Value = 44490
$ws.Range("J21").Value = 100

# Row 22
$ws.Range("D22").Value = 44399
$ws.Range("H22").Value = "Española"
$ws.Range("I22").Value = "Segunda"
$ws.Range("K22").Value = 15500
$ws.Range("L22").Value = 16000
$ws.Range("M22").Value = 15750
$ws.Range("P22").Value = 394

# Row 23
$ws.Range("D23").Value = 44491
$ws.Range("J23").Value = 100

# Row 24
$ws.Range("D24").Value = 44425
$ws.Range("K24").Value = 14000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 14500
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 362

# Row 25
$ws.Range("D25").Value = 44475
$ws.Range("K25").Value = 11000
$ws.Range("L25").Value = 12000
$ws.Range("M25").Value = 11500
$ws.Range("P25").Value = 288

# Row 26
$ws.Range("D26").Value = 44453
$ws.Range("J26").Value = 160
$ws.Range("K26").Value = 12500
$ws.Range("L26").Value = 13000
$ws.Range("M26").Value = 12750
$ws.Range("P26").Value = 319

# Row 27
$ws.Range("D27").Value = 44487
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 11000
$ws.Range("L27").Value = 12000
$ws.Range("M27").Value = 11500
$ws.Range("O27").Value = "Provincia del Elquí"
$ws.Range("P27").Value = 288

# Row 29
$ws.Range("D29").Value = 44467
$ws.Range("J29").Value = 160
$ws.Range("K29").Value = 11000
$ws.Range("L29").Value = 12000
$ws.Range("M29").Value = 11500
$ws.Range("O29").Value = "Provincia de Limarí"
$ws.Range("P29").Value = 288

# Row 30
$ws.Range("D30").Value = 44496
$ws.Range("J30").Value = 120

# Row 31
$ws.Range("D31").Value = 44482
$ws.Range("J31").Value = 120
$ws.Range("K31").Value = 11000
$ws.Range("L31").Value = 12000
$ws.Range("M31").Value = 11500
$ws.Range("P31").Value = 288

# Row 32
$ws.Range("D32").Value = 44455
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 13000
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = 13500
$ws.Range("P32").Value = 338

# Row 33
$ws.Range("D33").Value = 44432
$ws.Range("J33").Value = 120
$ws.Range("K33").Value = 14000
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = 14500
$ws.Range("P33").Value = 362

# Row 34
$ws.Range("D34").Value = 44515
$ws.Range("H34").Value = "Madrigal"
$ws.Range("I34").Value = "Primera"
$ws.Range("K34").Value = 11000
$ws.Range("L34").Value = 12000
$ws.Range("M34").Value = 11500
$ws.Range("P34").Value = 288

# Row 35
$ws.Range("D35").Value = 44508
$ws.Range("J35").Value = 160
$ws.Range("K35").Value = 11000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = 11500
$ws.Range("P35").Value = 288
